# Update the "想去人数" (want-to-go count) column F figures across the three
# sheets that carry this data (展览, 演出, 全部类型). 本地生活 has no data rows.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 471
$ws.Range("F5").Value = 75
$ws.Range("F6").Value = 14
$ws.Range("F7").Value = 1306
$ws.Range("F10").Value = 1312
$ws.Range("F13").Value = 25
$ws.Range("F16").Value = 107
$ws.Range("F21").Value = 232
$ws.Range("F22").Value = 2419
$ws.Range("F23").Value = 12
$ws.Range("F26").Value = 933
$ws.Range("F30").Value = 2828
$ws.Range("F34").Value = 676
$ws.Range("F36").Value = 1850
$ws.Range("F38").Value = 1863
$ws.Range("F39").Value = 204
$ws.Range("F42").Value = 42
$ws.Range("F43").Value = 876
$ws.Range("F44").Value = 801
$ws.Range("F45").Value = 1025
$ws.Range("F46").Value = 98
$ws.Range("F47").Value = 441
$ws.Range("F48").Value = 222
$ws.Range("F49").Value = 3348

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 10

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 471
$ws.Range("F5").Value = 75
$ws.Range("F8").Value = 14
$ws.Range("F9").Value = 1306
$ws.Range("F12").Value = 1312
$ws.Range("F15").Value = 25
$ws.Range("F18").Value = 107
$ws.Range("F23").Value = 232
$ws.Range("F24").Value = 2419
$ws.Range("F28").Value = 2828
$ws.Range("F34").Value = 676
$ws.Range("F36").Value = 1850
$ws.Range("F37").Value = 10
$ws.Range("F39").Value = 1863
$ws.Range("F41").Value = 876
$ws.Range("F42").Value = 801
$ws.Range("F43").Value = 1025
$ws.Range("F44").Value = 98
$ws.Range("F45").Value = 441
$ws.Range("F47").Value = 222
$ws.Range("F48").Value = 3348
